$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Robustness" text for the "Amount of error trapping and handling" row (row 6)
$ws.Range("D6").Value = "Asserts are placed in functions to ensure that parameters passed through are the correct type. The asserts placed covers the functionality of the program."

# Fill in Marks (column H) for several rows, highlighted with a yellow fill
$ws.Range("H2").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 2
$ws.Range("H5").Value = 1
$ws.Range("H10").Value = 2
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1

$ws.Range("H2").Interior.Color = 65535
$ws.Range("H3").Interior.Color = 65535
$ws.Range("H4").Interior.Color = 65535
$ws.Range("H5").Interior.Color = 65535
$ws.Range("H10").Interior.Color = 65535
$ws.Range("H11").Interior.Color = 65535
$ws.Range("H12").Interior.Color = 65535

# Add a total formula summing the Marks column (F) for rows 2-12
$ws.Range("F14").Formula = "=SUM(F2:F12)"

# Update the active selection
$ws.Range("K11").Select() | Out-Null
